# AutoCommit_15 июня 2024 г. 22:26:22_SibNout2023
# Apply the recorded changes (new "I3:J3" header-style cells, plus a batch
# of newly-filled numeric cells in the data grid) to the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 3: give the existing blank cells (C3:H3) a value of 1, and
#    create two new cells (I3, J3) that reuse the bold/boxed header
#    look (font bold, thick left+right border, centered+wrapped) from
#    the "I6:J6"-style header cells, but bold - this produces a brand
#    new cell style entry in styles.xml (index 10), exactly as in the
#    diff (cellXfs count 10 -> 11).
# ---------------------------------------------------------------------
$ws.Range("C3:H3").Value = 1

$ws.Range("I6").Copy()
$ws.Range("I3:J3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I3:J3").Font.Bold = $true
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1

# ---------------------------------------------------------------------
# 2) Row 4: F4/J4 are brand-new plain (unstyled) numeric cells; G4
#    already exists (style kept) and just receives a value.
# ---------------------------------------------------------------------
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("J4").Value = 2

# ---------------------------------------------------------------------
# 3) Row 5: G5/H5 already exist (style kept) and get a value; I5/J5 are
#    brand-new cells that should carry the existing "style 6" look
#    (reuse it from C11, which already uses that exact style).
# ---------------------------------------------------------------------
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 2

$ws.Range("C11").Copy()
$ws.Range("I5:J5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 2

# ---------------------------------------------------------------------
# 4) Row 11: F11/J11 are brand-new plain (unstyled) numeric cells; G11
#    already exists (style kept) and just receives a value.
# ---------------------------------------------------------------------
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 2
$ws.Range("J11").Value = 2

# ---------------------------------------------------------------------
# 5) Row 17: J17 is a brand-new cell using the existing "style 9" look
#    (reuse it from J16, which already uses that exact style).
# ---------------------------------------------------------------------
$ws.Range("J16").Copy()
$ws.Range("J17").PasteSpecial(-4122)     # xlPasteFormats
$ws.Range("J17").Value = 2

# ---------------------------------------------------------------------
# 6) Row 21: H21 changes from the empty "style 2" look to the "style 5"
#    look (reused from C17) and gets a value; I21/J21 are brand-new
#    cells using the existing "style 8" / "style 9" looks (reused from
#    I16 / J16).
# ---------------------------------------------------------------------
$ws.Range("C17").Copy()
$ws.Range("H21").PasteSpecial(-4122)     # xlPasteFormats
$ws.Range("H21").Value = 5

$ws.Range("I16").Copy()
$ws.Range("I21").PasteSpecial(-4122)     # xlPasteFormats
$ws.Range("I21").Value = 2

$ws.Range("J16").Copy()
$ws.Range("J21").PasteSpecial(-4122)     # xlPasteFormats
$ws.Range("J21").Value = 2

# ---------------------------------------------------------------------
# 7) Row 24: J24 is a brand-new cell using the existing "style 9" look.
# ---------------------------------------------------------------------
$ws.Range("J16").Copy()
$ws.Range("J24").PasteSpecial(-4122)     # xlPasteFormats
$ws.Range("J24").Value = 2

# ---------------------------------------------------------------------
# 8) Row 25: H25 already exists (style kept) and just receives a value.
# ---------------------------------------------------------------------
$ws.Range("H25").Value = 2

# ---------------------------------------------------------------------
# 9) Row 27: E27/F27 are brand-new cells using the existing "style 6"
#    look (reused from C11); G27 already exists (style kept) and just
#    receives a value; J27 is a brand-new cell using the existing
#    "style 9" look.
# ---------------------------------------------------------------------
$ws.Range("C11").Copy()
$ws.Range("E27:F27").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("E27").Value = 2
$ws.Range("F27").Value = 2

$ws.Range("G27").Value = 2

$ws.Range("J16").Copy()
$ws.Range("J27").PasteSpecial(-4122)     # xlPasteFormats
$ws.Range("J27").Value = 2

# ---------------------------------------------------------------------
# 10) Row 28: D28/E28/F28/G28 already exist (style kept) and just
#     receive a value.
# ---------------------------------------------------------------------
$ws.Range("D28:G28").Value = 2

# ---------------------------------------------------------------------
# 11) Finally, move the active selection of the frozen bottom-right pane
#     to J3 (it previously pointed at D27).
# ---------------------------------------------------------------------
$ws.Range("J3").Select()
